$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    $origStyle = $range.Style
    $range.NumberFormat = "@"
    $range.Value = $value
    $range.Style = $origStyle
}

Set-TextValue $ws.Range("D2") "26.514.02"
Set-TextValue $ws.Range("E2") "  +0.41%  "
Set-TextValue $ws.Range("D3") "1.630.50"
Set-TextValue $ws.Range("D4") "0.998"
Set-TextValue $ws.Range("E4") "  -0.13%  "
Set-TextValue $ws.Range("D5") "213.17"
Set-TextValue $ws.Range("E5") "  -0.31%  "
Set-TextValue $ws.Range("D6") "0.503"
Set-TextValue $ws.Range("E6") "  -0.06%  "
Set-TextValue $ws.Range("D7") "0.999"
Set-TextValue $ws.Range("E7") "  -0.14%  "
Set-TextValue $ws.Range("E8") "  +0.32%  "
Set-TextValue $ws.Range("E9") "  +0.16%  "
Set-TextValue $ws.Range("E10") "  +0.23%  "
Set-TextValue $ws.Range("D11") "0.0851"
Set-TextValue $ws.Range("E11") "  -0.44%  "
Set-TextValue $ws.Range("D12") "1.860.28"
Set-TextValue $ws.Range("E12") "  +0.99%  "
Set-TextValue $ws.Range("D13") "1.627.86"
Set-TextValue $ws.Range("E13") "  +0.58%  "
Set-TextValue $ws.Range("E14") "  +0.05%  "
Set-TextValue $ws.Range("E15") "  +0.27%  "
Set-TextValue $ws.Range("D16") "64.08"
Set-TextValue $ws.Range("E16") "  -0.87%  "
Set-TextValue $ws.Range("D17") "237.29"
Set-TextValue $ws.Range("E17") "  +4.33%  "
Set-TextValue $ws.Range("D18") "26.531.86"
Set-TextValue $ws.Range("E18") "  +0.40%  "
Set-TextValue $ws.Range("D19") "7.83"
Set-TextValue $ws.Range("E19") "  +3.83%  "
Set-TextValue $ws.Range("D20") "0.0₃0727"
Set-TextValue $ws.Range("E20") "  +0.14%  "
Set-TextValue $ws.Range("E21") "  -0.09%  "
Set-TextValue $ws.Range("E22") "  -0.86%  "
Set-TextValue $ws.Range("E23") "  +0.55%  "
Set-TextValue $ws.Range("E24") "  +2.49%  "
Set-TextValue $ws.Range("D25") "147.01"
Set-TextValue $ws.Range("E25") "  +1.29%  "
Set-TextValue $ws.Range("E26") "  -0.14%  "
Set-TextValue $ws.Range("E27") "  +0.97%  "
Set-TextValue $ws.Range("E28") "  +0.04%  "
Set-TextValue $ws.Range("D29") "15.70"
Set-TextValue $ws.Range("E29") "  +2.14%  "
Set-TextValue $ws.Range("E30") "  +0.14%  "
Set-TextValue $ws.Range("E31") "  -0.35%  "
Set-TextValue $ws.Range("D32") "1.522.65"
Set-TextValue $ws.Range("E32") "  +5.05%  "
Set-TextValue $ws.Range("E33") "  +1.23%  "
Set-TextValue $ws.Range("E34") "  -0.22%  "
Set-TextValue $ws.Range("D35") "1.52"
Set-TextValue $ws.Range("E35") "  +3.23%  "
Set-TextValue $ws.Range("E36") "  -0.07%  "
Set-TextValue $ws.Range("E37") "  +2.13%  "
Set-TextValue $ws.Range("E38") "  +0.13%  "
Set-TextValue $ws.Range("E39") "  +0.23%  "
Set-TextValue $ws.Range("D40") "5.88"
Set-TextValue $ws.Range("E40") "  +0.47%  "
Set-TextValue $ws.Range("D41") "0.998"
Set-TextValue $ws.Range("E41") "  -0.12%  "
Set-TextValue $ws.Range("E42") "  +0.70%  "
Set-TextValue $ws.Range("D43") "1.770.80"
Set-TextValue $ws.Range("E43") "  +0.94%  "
Set-TextValue $ws.Range("D44") "63.10"
Set-TextValue $ws.Range("E44") "  +1.88%  "
Set-TextValue $ws.Range("E45") "  -0.11%  "
Set-TextValue $ws.Range("D46") "0.907"
Set-TextValue $ws.Range("E46") "  -0.76%  "
Set-TextValue $ws.Range("D47") "90.41"
Set-TextValue $ws.Range("E47") "  +2.89%  "
Set-TextValue $ws.Range("E48") "  +1.92%  "
Set-TextValue $ws.Range("E49") "  -3.46%  "
Set-TextValue $ws.Range("D50") "0.0501"
Set-TextValue $ws.Range("E50") "  -0.25%  "
Set-TextValue $ws.Range("D51") "0.0966"
Set-TextValue $ws.Range("E51") "  +0.58%  "
